# The sheet previously held two rows of "class" example data in A1:E2.
# Split/rebalance it into three rows of "class + gender" example data
# spanning A1:F3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (6 entries)
$ws.Range("A1").Value = "예시111"
$ws.Range("B1").Value = "예시112"
$ws.Range("C1").Value = "예시113"
$ws.Range("D1").Value = "예시211"
$ws.Range("E1").Value = "예시213"
$ws.Range("F1").Value = "예시314"

# Row 2 (2 entries)
$ws.Range("A2").Value = "예시510"
$ws.Range("B2").Value = "예시415"

# Row 3 (2 entries)
$ws.Range("A3").Value = "예시401"
$ws.Range("B3").Value = "예시619"

# Newly-used cells (F1, A3, B3) fall outside the original A1:E2 range and
# don't inherit its formatting automatically, so copy the existing cell
# style (from A1) onto them without touching any other, still-blank cells.
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Page setup (paper size / orientation) for printing
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("A4").Select()
